# Weekly update: insert a new price record for "Arveja Verde" at
# Vega Modelo de Temuco, shifting the existing rows 28-46 down to 29-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 28 (pushes rows 28..46 down to 29..47)
$ws.Rows.Item(28).EntireRow.Insert()

# Populate the newly inserted row 28 with the new record's data
$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 44484
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 100112022
$ws.Range("G28").Value = "Arveja Verde"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("M28").Value = 25000
$ws.Range("N28").Value = "`$/malla 25 kilos"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 1000
$ws.Range("Q28").Value = 25
$ws.Range("R28").Value = "Hortaliza"
